# Generate Report for Handoff
#
# The handoff XLIFF files were (re)generated, so the localization-status
# report needs to reflect the new status and the refreshed timestamps:
#   - Status moves from "In Translation" to "Ready for handoff"
#   - The "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#     timestamps are bumped to the new generation time
#   - The Status column is widened so the new (longer) text fits

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Refreshed handoff-generation timestamps ---
$overview.Range("G2").Value = "2016-08-31 13:14:21"
$dede.Range("H2").Value     = "2016-08-31 13:14:21"
$zhcn.Range("H2").Value     = "2016-08-31 13:14:17"

# --- Widen the Status columns to fit "Ready for handoff" ---
$overview.Range("E:E").ColumnWidth = 16.33333333
$overview.Range("F:F").ColumnWidth = 16.33333333
$zhcn.Range("C:C").ColumnWidth     = 16.33333333
$dede.Range("C:C").ColumnWidth     = 16.33333333
